# Refresh cryptocurrency Price / Volume(1h) columns on Sheet1 with the
# latest scraped figures (GitHub Actions scheduled update), and fix the
# row order for TheSandbox / EnergySwap which swapped rank positions.
#
# The sheet stores Price/Volume as text (not numbers) so formatting such
# as "29.227.67" / "  +0.28%  " is preserved exactly. Values that look
# like plain decimals (e.g. "0.9991") are entered with a leading
# apostrophe so Excel keeps them as text instead of auto-converting them
# to numbers; .Style is then reset to "Normal" so the quote-prefix
# indicator doesn't leave a stray cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.227.67"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.833.89"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'243.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.6202"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").Value = "'0.2898"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "'23.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "'0.07677"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "1.824.89"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'4.981"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "'0.6709"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'82.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'0.000008969"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").Value = "'5.875"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "29.202.16"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "2.070.03"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("D20").Value = "'236.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'7.356"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").Value = "'0.9989"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "'158.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "'0.1401"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").Value = "'8.563"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "'1.491"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").Value = "'0.05767"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("D31").Value = "'4.111"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'4.093"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "'1.208"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "'1.871"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "'0.7323"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("D36").Value = "'1.144"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "'2.605"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("D38").Value = "'2.858"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "1.226.97"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").Value = "'6.261"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").Value = "'0.9098"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "'101.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "1.974.49"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "'65.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'0.5038"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.171"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4030"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'0.1137"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.45%  "
